$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.266.21"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "1.592.91"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'213.07"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("E10").Value = "  -2.07%  "
$ws.Range("D11").Value = "'0.0851"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "1.817.85"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "1.566.52"
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("E15").Value = "  -2.36%  "
$ws.Range("D16").Value = "'63.77"
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("D17").Value = "26.262.69"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").Value = "'215.45"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").Value = "'7.35"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "'4.28"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "'9.02"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  -2.75%  "
$ws.Range("D25").Value = "'144.80"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'6.95"
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("D29").Value = "'15.11"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("D30").Value = "'0.0489"
$ws.Range("E30").Value = "  -2.30%  "
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").Value = "1.411.39"
$ws.Range("E33").Value = "  +5.47%  "
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("D37").Value = "'0.571"
$ws.Range("E37").Value = "  -5.27%  "
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("D39").Value = "'0.822"
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("D40").Value = "'5.77"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "'0.948"
$ws.Range("E42").Value = "  -9.98%  "
$ws.Range("E43").Value = "  +0.93%  "
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").Value = "1.729.69"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "'60.89"
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("D47").Value = "'87.20"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("E48").Value = "  -1.39%  "
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("D50").Value = "'0.0952"
$ws.Range("E50").Value = "  -3.02%  "
$ws.Range("E51").Value = "  +0.08%  "
